# feat: add 2022-Q3 data
#
# 1. Insert a new "2022-Q3" row at the top of the summary ("总计") sheet,
#    pushing the existing quarters down and renumbering the index column.
# 2. Insert a brand-new "2022-Q3" worksheet (right after "总计") holding the
#    per-fund holdings detail for that quarter, using the same layout/style
#    as the existing quarter sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Part 1: update the "总计" (summary) sheet
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)

# Push rows 2..4 down to 3..5, leaving a blank row 2 for the new quarter.
$summary.Rows.Item(2).Insert()

# The freshly-inserted row inherited row 1's (header) formatting; strip it
# back to a plain cell, then restore the bold/centered index-column style
# (style "2") on A2 by copying it from A3 (the old A2, now shifted down).
$summary.Range("A2:D2").ClearFormats()
$summary.Range("A3").Copy()
$summary.Range("A2").PasteSpecial(-4122)

$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q3"
$summary.Range("C2").Value = 17
$summary.Range("D2").Value = 1.47

# Renumber the index column so it stays a contiguous 0..3 sequence.
$summary.Range("A3").Value = 1
$summary.Range("A4").Value = 2
$summary.Range("A5").Value = 3

# ---------------------------------------------------------------------------
# Part 2: add the new "2022-Q3" detail worksheet
# ---------------------------------------------------------------------------
$newSheet = $wb.Worksheets.Add($null, $summary)
$newSheet.Name = "2022-Q3"

# Borrow the header-row and data-row formatting from the existing "2022-Q1"
# detail sheet (same column layout: B..H headers, A-column style "2").
$template = $wb.Worksheets.Item("2022-Q1")
$template.Range("A1:H1").Copy($newSheet.Range("A1:H1"))
$template.Range("A2:H2").Copy($newSheet.Range("A2:H18"))

# Header row.
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Per-fund holdings data. Columns B, D, E, F (and G except the last row) are
# numeric-looking text in the source workbook (e.g. "0.50", "011136") so we
# must force a text number-format before assigning them, otherwise Excel's
# COM layer auto-coerces them into numbers and mangles the formatting
# (leading zeros / trailing zeros get dropped).
$data = @(
    @('011136', '广发盛兴混合A', '16.84', '94.77', '3.37', '0.5675', 10),
    @('506007', '广发科创板两年定开混合', '5.01', '94.25', '6.59', '0.3302', 2),
    @('012342', '广发瑞泽精选混合A', '4.99', '93.90', '3.23', '0.1612', 10),
    @('009874', '九泰久睿量化股票A', '3.15', '93.83', '3.41', '0.1074', 5),
    @('013000', '广发盛泽一年持有期混合A', '2.29', '82.39', '3.34', '0.0765', 10),
    @('002133', '广发鑫益灵活配置混合', '1.22', '93.85', '5.86', '0.0715', 5),
    @('011137', '广发盛兴混合C', '1.74', '94.77', '3.37', '0.0586', 10),
    @('010120', '九泰久福量化股票A', '0.54', '93.91', '3.44', '0.0186', 4),
    @('001897', '九泰久盛量化先锋灵活配置混合A', '0.50', '93.59', '3.40', '0.0170', 5),
    @('011500', '九泰量化新兴产业混合', '0.58', '93.85', '2.81', '0.0163', 7),
    @('009043', '九泰久信量化股票', '0.43', '93.60', '3.42', '0.0147', 5),
    @('012343', '广发瑞泽精选混合C', '0.36', '93.90', '3.23', '0.0116', 10),
    @('004510', '九泰久盛量化先锋灵活配置混合C', '0.28', '93.59', '3.40', '0.0095', 5),
    @('013001', '广发盛泽一年持有期混合C', '0.27', '82.39', '3.34', '0.0090', 10),
    @('005360', '汇安资产轮动灵活配置混合', '0.12', '80.44', '2.95', '0.0035', 10),
    @('010121', '九泰久福量化股票C', '0.04', '93.91', '3.44', '0.0014', 4),
    @('016399', '九泰久睿量化股票C', '0.00', '93.83', '3.41', $null, 5)
)

$row = 2
foreach ($fund in $data) {
    $newSheet.Range("A$row").Value = $row - 2

    $textRange = $newSheet.Range("B$row:F$row")
    $textRange.NumberFormat = "@"
    $newSheet.Range("B$row").Value = $fund[0]
    $newSheet.Range("C$row").Value = $fund[1]
    $newSheet.Range("D$row").Value = $fund[2]
    $newSheet.Range("E$row").Value = $fund[3]
    $newSheet.Range("F$row").Value = $fund[4]
    $textRange.Style = "Normal"

    if ($null -eq $fund[5]) {
        $newSheet.Range("G$row").Value = 0
    } else {
        $newSheet.Range("G$row").NumberFormat = "@"
        $newSheet.Range("G$row").Value = $fund[5]
        $newSheet.Range("G$row").Style = "Normal"
    }

    $newSheet.Range("H$row").Value = $fund[6]

    $row++
}
